$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header (H1) - copy formatting from the neighboring
# header cell (G1) so it matches the existing header style, then set text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column (H2)
$ws.Range("H2").Value = 0
